$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 26290626

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 14.9
